# Powerpoint writer: consolidate text run nodes.
#
# Slide 1, shape 2 ("TextBox 3") holds the caption as five separate
# runs: "The" / " " / "picture" / " " / "first". Fold each run that
# is immediately followed by a lone-space run into that run (so
# "The"+" " -> "The " and "picture"+" " -> "picture "), leaving
# "first" untouched.
#
# Characters(...).InsertAfter() appends into the run that exactly
# covers that character range instead of minting a fresh run, so
# copying the adjacent space onto the end of the preceding run (then
# deleting the now-redundant standalone space run) merges the two
# without disturbing the kept run's existing (empty) <a:rPr/>.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# "The" + " " -> "The "
$run1 = $tr.Characters(1, 3)
[void]$run1.InsertAfter(" ")
[void]$tr.Characters(5, 1).Delete()

# Text is now "The picture first"; "picture" + " " -> "picture "
$run2 = $tr.Characters(5, 7)
[void]$run2.InsertAfter(" ")
[void]$tr.Characters(13, 1).Delete()
